$d = $word.ActiveDocument

# 1. Add "Multivariate Statistics" text into the empty Heading3 paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 3") {
        $r = $p.Range
        $r.Collapse(1) | Out-Null
        $r.Text = "Multivariate Statistics"
        break
    }
}

# 2. Ativação date change.
# NB: a plain Find/Replace (or a direct Range.Text= on just this run) on this
# particular run causes the COM host to coalesce it with the following
# "Departamento: ..." run (losing the <w:br/> run boundary between them).
# Rebuild the "Ativação: .../Departamento: ..." span explicitly via
# InsertAfter so the original two-run/<w:br/> structure is preserved.
$rngAtiv = $d.Content
$rngAtiv.Find.Execute("Ativação: 01/01/2016") | Out-Null
$ativStart = $rngAtiv.Start
$rngDepto = $d.Content
$rngDepto.Find.Execute("Departamento: Ciências Básicas e Ambientais") | Out-Null
$ativEnd = $rngDepto.End
$rngSpan = $d.Range($ativStart, $ativEnd)
$rngSpan.Text = ""
$rngSpan.Collapse(0) | Out-Null
$rngSpan.InsertAfter("Ativação: 01/01/2021" + [char]11)
$rngSpan.Collapse(0) | Out-Null
$rngSpan.InsertAfter("Departamento: Ciências Básicas e Ambientais")

# 3. Programa resumido paragraph replace.
$old3 = "Distribuições Multivariadas, Confiabilidade, Distribuições Amostrais Multivariadas, Inferências Sobre Vetores de Médias, Testes de Hipóteses Sobre Médias, Associação entre variáveis , Regressão Logística, Análise de Variância Multivariada, Análise de Agrupamento, Componentes Principais, Análise Fatorial"
$new3 = "Probabilidade: Vetor de variáveis aleatórias, Distribuição conjunta/marginal, Esperança e variância condicional/marginal. Estatística: Regressão Logística simples, Teste Qui-Quadrado, Testes de normalidade, Testes não-paramétricos. Técnicas Multivariadas: Gráficos multivariados, Regressão Linear Múltipla, Regressão Logística Múltipla, Análise de Variância Múltipla; Análise de agrupamento; Análise de componentes principais; Análise fatorial; Análise discriminante e Análise de correspondência"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new3, 2) | Out-Null

# 4. Programa paragraph - collapse multi-run/break content to a single run.
$newProgramaText = "Probabilidade: Vetor de variáveis aleatórias, Distribuição conjunta/marginal, Esperança e Variância condicional/marginal.Estatística: Regressão Logística simples (coeficiente de associação, sensitividade e especificidade, risco relativo, razão de chances), Teste Qui-Quadrado (testes de aderência, homogeneidade e independência), Teste de normalidade (Shapiro-Wilk, Teste de Kolmogorov-Smirnov), Testes não-paramétricos para amostras pareadas e independentes.Técnicas Multivariadas: Gráficos multivariados, Regressão Linear Múltipla, Regressão Logística Múltipla, Análise de Variância Múltipla; Análise de agrupamento; Análise de componentes principais; Análise fatorial; Análise discriminante e Análise de correspondência"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Vetor de variáveis aleatórias, Distribuição conjunta, distribuição marginal")) {
        $pr = $p.Range
        $pr.MoveEnd(1, -1) | Out-Null
        $pr.Text = $newProgramaText
        break
    }
}

# 5. Método value replace.
$d.Content.Find.Execute("Aulas expositivas teóricas, aulas práticas, aulas de exercícios.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: N = (N1+...+Nn)/n", 2) | Out-Null

# 6. Critério value replace.
$d.Content.Find.Execute("Duas Avaliações P1 e P2 sendo a média calculada por (P1+2*P2)/3", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "NF≥ 5,0.", 2) | Out-Null

# 7. Norma de recuperação value replace.
$d.Content.Find.Execute("Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.", 2) | Out-Null

# 8. Bibliografia paragraph - collapse multi-run/break content to a single run.
$newBiblioText = "G.C. Runger, D. Montgomery. Estatística aplicada e probabilidade para engenheiros. São Paulo: Ed. LTC, 2009. D. C. Montgomery, E. A. Peck, G. G. Vining, Introduction to Linear Regression Analysis, 4th ed., Hoboken: John Wiley, 2006.W. J. Conover, Practical Nonparametric Statistics, 3rd ed., New York: John Wiley d Sons, 1999.R. A. Johnson, D. W. Wichern, Applied Multivariate Statistical Analysis, 6th ed., New Jersey: Prentice Hall, 2007."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("W. O. Bussab, P. A. Morettin")) {
        $pr = $p.Range
        $pr.MoveEnd(1, -1) | Out-Null
        $pr.Text = $newBiblioText
        break
    }
}
